# "updated legacy GSC export data"
# The GSC export rolled forward by one day: the oldest date (row 2,
# 2025-09-04) drops off, every remaining row's data shifts up one row,
# and two freshly-observed days are appended at the bottom with zero
# counts (no data collected for them yet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$lastRow = $ws.Cells.Item(1, 1).End(4).Row   # 90 (header is row 1)

# Read the existing date labels / HTTPS URL counts for rows 3..lastRow
# (i.e. everything except the oldest row, which is being dropped), then
# tack on the two new trailing days.
$dates = @()
$httpsUrls = @()
for ($r = 3; $r -le $lastRow; $r++) {
    $dates += $ws.Cells.Item($r, 1).Value()
    $httpsUrls += $ws.Cells.Item($r, 3).Value()
}
$dates += "2025-12-02"
$dates += "2025-12-03"
$httpsUrls += 0
$httpsUrls += 0

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.Value = "'" + $dates[$i]
    $dateCell.ClearFormats()
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = $httpsUrls[$i]
}
